# Implemented password masking using Base64 encoding & decoding
#
# This script edits testdata.xlsx so that plaintext passwords stored on the
# "DataProviderTests" sheet are replaced with their Base64-encoded form
# (admin123 / admin12345 -> YWRtaW4xMjM=, the Base64 encoding of "admin123"),
# flips a few "execute" yes/no flags, restyles the header/body of the
# "Tests" sheet (left-aligned instead of center-aligned), widens column E on
# "DataProviderTests" to fit the new values, and moves the active
# selection/tab the way the author last left the workbook.

$wb = $excel.ActiveWorkbook
$wsTests = $wb.Worksheets.Item(1)          # "Tests"
$wsData  = $wb.Worksheets.Item(2)          # "DataProviderTests"

# ---------------------------------------------------------------------
# 1) Tests sheet: restyle header (row 1) and body (rows 2-6) of columns
#    A & B from center-aligned to left-aligned, and flip two execute flags.
# ---------------------------------------------------------------------
$wsTests.Range("A1:B1").HorizontalAlignment = -4131   # xlLeft
$wsTests.Range("A1:B1").VerticalAlignment   = -4108   # xlCenter

$wsTests.Range("A2:B6").HorizontalAlignment = -4131   # xlLeft
$wsTests.Range("A2:B6").VerticalAlignment   = -4108   # xlCenter

$wsTests.Range("C2").Value = "yes"
$wsTests.Range("C6").Value = "no"

# ---------------------------------------------------------------------
# 2) DataProviderTests sheet: mask passwords with Base64, flip a few
#    execute flags, and widen column E to fit the new values.
# ---------------------------------------------------------------------
$wsData.Range("E2").Value  = "YWRtaW4xMjM="
$wsData.Range("E3").Value  = "YWRtaW4xMjM="
$wsData.Range("E4").Value  = "YWRtaW4xMjM="
$wsData.Range("E5").Value  = "YWRtaW4xMjM="
$wsData.Range("E6").Value  = "YWRtaW4xMjM="
$wsData.Range("E7").Value  = "YWRtaW4xMjM="
$wsData.Range("E8").Value  = "YWRtaW4xMjM="
$wsData.Range("E9").Value  = "YWRtaW4xMjM="

$wsData.Range("B2").Value  = "yes"
$wsData.Range("B3").Value  = "yes"
$wsData.Range("B10").Value = "no"

$wsData.Columns.Item(5).ColumnWidth = 14.8

# ---------------------------------------------------------------------
# 3) Selection / active tab: the author left the workbook with
#    "Tests"!C2 selected and "DataProviderTests" as the active tab with
#    E15 selected.
# ---------------------------------------------------------------------
$wsTests.Activate()
$wsTests.Range("C2").Select()

$wsData.Activate()
$wsData.Range("E15").Select()

"done"
